$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Text constants (kept exact, including curly apostrophes / trailing
# spaces, so they reuse the existing shared-string entries).
# ------------------------------------------------------------------
$ruleText1 = "All must play 3 innings in infield, some may play 4 but no more. "
$ruleText2 = "Possibly have one kid sit on bench each inning if full team shows up. "
$ruleText3 = "There’s also a few kids i don’t want to play certain positions. If you can do that let me know. Mainly 1st and PH cause they might get hurt. "

$yellow = 65535

# ------------------------------------------------------------------
# Step 1: add the new "Details" sheet right after "Rules", and copy
# the long-form rule text into it *before* we touch the Rules sheet,
# so the shared-string entries survive the upcoming edits.
# ------------------------------------------------------------------
$rules = $wb.Worksheets.Item("Rules")
$details = $wb.Worksheets.Add($null, $rules)
$details.Name = "Details"

$details.Range("A1").Value = $ruleText1
$details.Range("A2").Value = $ruleText2
$details.Range("A3").Value = $ruleText3

$details.Columns.Item(1).ColumnWidth = 124.59244791666667
$details.Range("A6").Select()

# ------------------------------------------------------------------
# Step 2: rebuild the "Rules" sheet content into the new rule-table
# layout, then rename it.
# ------------------------------------------------------------------
$rules.Range("B1").Value = "Rule_Word"
$rules.Range("C1").Value = "Value"

$rules.Range("A2").Value = "All"
$rules.Range("B2").Value = "Min_Infield_Inning"
$rules.Range("C2").Value = 3

$rules.Range("A3").Value = "Grayson"
$rules.Range("B3").Value = "Position_Restricted"
$rules.Range("C3").Value = "1st"

$rules.Range("A4").ClearContents()

$rules.Columns.Item(1).ColumnWidth = 124.59244791666667
$rules.Columns.Item(2).ColumnWidth = 17.736979166666668
$rules.Columns.Item(3).ColumnWidth = 5.307291666666667

$rules.Range("A19").Select()
$rules.Name = "Rules - Not working"

# ------------------------------------------------------------------
# Step 3: highlight the Players sheet (everyone except "Luke M" in
# row 7) and move the selection.
# ------------------------------------------------------------------
$players = $wb.Worksheets.Item("Players")
$players.Range("A2:A6").Interior.Color = $yellow
$players.Range("A8:A12").Interior.Color = $yellow
$players.Range("A7").Select()

# ------------------------------------------------------------------
# Step 4: highlight the Positions sheet (every data row) and move the
# selection.
# ------------------------------------------------------------------
$positions = $wb.Worksheets.Item("Positions")
$positions.Range("A2:B11").Interior.Color = $yellow
$positions.Range("A5:B5").Select()

# ------------------------------------------------------------------
# Step 5: make "Rules - Not working" the active tab.
# ------------------------------------------------------------------
$rules.Activate()
